# Adds std/min/max "sucesso" metrics for each origem-modalidade group
# (aon, flex, sub) to the 2023 summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 3 new columns after each "*_media_sucesso" column -------
# Insert order matters: always insert right after "aon_media_sucesso"
# (column I) first, then after "flex_media_sucesso" (which, after the
# first insert, sits at column S), then after "sub_media_sucesso"
# (which, after the first two inserts, sits at column AC).

$ws.Range("J1:L1").EntireColumn.Insert()
$ws.Range("T1:V1").EntireColumn.Insert()
$ws.Range("AD1:AF1").EntireColumn.Insert()

# --- 2. Header row (row 1) ---------------------------------------------
$ws.Range("J1").Value  = "aon_std_sucesso"
$ws.Range("K1").Value  = "aon_min_sucesso"
$ws.Range("L1").Value  = "aon_max_sucesso"

$ws.Range("T1").Value  = "flex_std_sucesso"
$ws.Range("U1").Value  = "flex_min_sucesso"
$ws.Range("V1").Value  = "flex_max_sucesso"

$ws.Range("AD1").Value = "sub_std_sucesso"
$ws.Range("AE1").Value = "sub_min_sucesso"
$ws.Range("AF1").Value = "sub_max_sucesso"

# --- 3. Data rows (row 2 = apoia.se, row 3 = catarse) -------------------
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

$ws.Range("J3").Value = 44961.93536949201
$ws.Range("K3").Value = 41.81688448509265
$ws.Range("L3").Value = 679297.6600721752

$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0

$ws.Range("T3").Value = 33934.82811955066
$ws.Range("U3").Value = 10.77163914429046
$ws.Range("V3").Value = 708972.7845446636

$ws.Range("AD2").Value = 682.4025885496077
$ws.Range("AE2").Value = 1.087396962410123
$ws.Range("AF2").Value = 5087.076865717208

$ws.Range("AD3").Value = 198.3989605548985
$ws.Range("AE3").Value = 10.98162164796783
$ws.Range("AF3").Value = 538.4389998789497

# --- 4. Styling ----------------------------------------------------------
# Inserting the columns already copies the formatting of the column to
# their left, so the "std"/"min" data cells (J,K / T,U / AD,AE) come out
# with the same currency style (style index 2, "R$ #,##0.00") as the
# neighboring "*_valor_sucesso"/"*_media_sucesso" columns, and the header
# cells come out with the bold header style - exactly what's needed.
# The "max" columns (L / V / AF), however, must end up with the plain
# default (General) style on the data rows, so clear their style there.

$ws.Range("L2").Style = "Normal"
$ws.Range("V2").Style = "Normal"
$ws.Range("AF2").Style = "Normal"

$ws.Range("L3").Style = "Normal"
$ws.Range("V3").Style = "Normal"
$ws.Range("AF3").Style = "Normal"
